# [Snehil] Add : Automating Email using Blue Prism
#
# Update the two employee email addresses on Sheet1 and move the active
# selection to B3 (matches the committed workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Snehil's email address changes.
$ws.Range("B2").Value = "snehil123488@gmail.com"

# swastik's email address changes.
$ws.Range("B3").Value = "swastik123456@gmail.com"

# The active cell/selection moves from C3 to B3.
$ws.Range("B3").Select()
